# Fruta / hortaliza, semanal
# Inserts two new daily price rows (Flame Seedless and Ralli Seedless lots
# dated 2023-01-26) at the top of the data block (rows 39-40), pushing the
# existing historical rows down by two (old row 39 -> new row 41, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 39; Excel shifts rows 39:133
# down to 41:135 and extends the used range / dimension accordingly.
$ws.Rows("39:40").Insert()

# --- New row 39 ---------------------------------------------------------
$ws.Range("A39").Value = 8
$ws.Range("B39").Value = "Terminal La Palmera de La Serena"
$ws.Range("C39").Value = "Coquimbo"
$ws.Range("D39").Value = 44952
$ws.Range("E39").Value = 4
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100109
$ws.Range("H39").Value = "Uva"
$ws.Range("I39").Value = 100109001
$ws.Range("J39").Value = "Uva"
$ws.Range("K39").Value = "Flame Seedless"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 700
$ws.Range("N39").Value = 7000
$ws.Range("O39").Value = 8000
$ws.Range("P39").Value = 7500
$ws.Range("Q39").Value = "`$/bandeja 18 kilos"
$ws.Range("R39").Value = "Provincia del Elquí"
$ws.Range("S39").Value = 417
$ws.Range("T39").Value = 18

# --- New row 40 ---------------------------------------------------------
$ws.Range("A40").Value = 8
$ws.Range("B40").Value = "Terminal La Palmera de La Serena"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 44952
$ws.Range("E40").Value = 4
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100109
$ws.Range("H40").Value = "Uva"
$ws.Range("I40").Value = 100109001
$ws.Range("J40").Value = "Uva"
$ws.Range("K40").Value = "Ralli Seedless"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 400
$ws.Range("N40").Value = 9000
$ws.Range("O40").Value = 10000
$ws.Range("P40").Value = 9500
$ws.Range("Q40").Value = "`$/bandeja 18 kilos"
$ws.Range("R40").Value = "Provincia del Elquí"
$ws.Range("S40").Value = 528
$ws.Range("T40").Value = 18
